$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# --- Header row: "request parameter" -> "request_parameter" ---
$ws.Range("F1").Value = "request_parameter"

# --- Row 2: success case ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "success_case "
$ws.Range("C2").Value = "post"
$ws.Range("D2").Value = "https://openapiv5.ketangpai.com/UserApi/login"
$ws.Range("E2").Value = "{""Content-Type"":""application/json""}"
$ws.Range("F2").Value = "{""email"":""2378807139@qq.com "",""password"":""123456""}"
$ws.Range("G2").Value = "{""status"":1,""code"":10000,""message"":""访问成功""}"

# --- Row 3: failed case ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "failed_case"
$ws.Range("C3").Value = "post"
$ws.Range("D3").Value = "https://openapiv5.ketangpai.com/UserApi/login"
$ws.Range("E3").Value = "{""Content-Type"":""application/json""}"
$ws.Range("F3").Value = "{""email"":""2378807139@qq.com "",""password"":""lemon""}"
$ws.Range("G3").Value = "{""status"":0,""code"":""30508"",""message"":""登录失败""}"

# --- Hyperlinks on D2 / D3 ---
$ws.Hyperlinks.Add($ws.Range("D2"), "https://openapiv5.ketangpai.com/UserApi/login")
$ws.Hyperlinks.Add($ws.Range("D3"), "https://openapiv5.ketangpai.com/UserApi/login")

# --- Selection moves to D6 ---
$ws.Range("D6").Select()
